$wb = $excel.ActiveWorkbook

# =========================================================
# Step 1: insert a new "2022-Q3" sheet immediately before the
# existing "2022-Q2" sheet (all later quarters shift right).
# =========================================================
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"
# Re-fetch by name -- the object handed back by Add() can be stale
# once the sheet collection is re-indexed, so always look sheets up
# by name before touching them further.
$ws2 = $wb.Worksheets.Item("2022-Q3")

# --- header row: clone formatting (bold/border) from a sibling sheet, then set text ---
$hdrSrc = $wb.Worksheets.Item("2022-Q2").Range("B1:H1")
$hdrSrc.Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws2.Cells.Item(1,2).Value = "基金代码"
$ws2.Cells.Item(1,3).Value = "基金名称"
$ws2.Cells.Item(1,4).Value = "基金规模"
$ws2.Cells.Item(1,5).Value = "股票总仓位"
$ws2.Cells.Item(1,6).Value = "仓位占比"
$ws2.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1,8).Value = "仓位排名"

# --- column A (row index) formatting for all 40 data rows ---
$aSrc = $wb.Worksheets.Item("2022-Q2").Range("A2")
$aSrc.Copy()
$ws2.Range("A2:A41").PasteSpecial(-4122)

# --- data rows: A (index, number), B..G (text), H (rank, number) ---
$ws2.Cells.Item(2,1).Value = 0
$c = $ws2.Cells.Item(2,2); $c.NumberFormat = "@"; $c.Value = "005299"
$ws2.Cells.Item(2,3).Value = "万家成长优选灵活配置混合A"
$c = $ws2.Cells.Item(2,4); $c.NumberFormat = "@"; $c.Value = "14.52"
$c = $ws2.Cells.Item(2,5); $c.NumberFormat = "@"; $c.Value = "94.06"
$c = $ws2.Cells.Item(2,6); $c.NumberFormat = "@"; $c.Value = "2.94"
$c = $ws2.Cells.Item(2,7); $c.NumberFormat = "@"; $c.Value = "0.4269"
$ws2.Cells.Item(2,8).Value = 10
$ws2.Cells.Item(3,1).Value = 1
$c = $ws2.Cells.Item(3,2); $c.NumberFormat = "@"; $c.Value = "003751"
$ws2.Cells.Item(3,3).Value = "万家瑞隆混合A"
$c = $ws2.Cells.Item(3,4); $c.NumberFormat = "@"; $c.Value = "14.54"
$c = $ws2.Cells.Item(3,5); $c.NumberFormat = "@"; $c.Value = "93.96"
$c = $ws2.Cells.Item(3,6); $c.NumberFormat = "@"; $c.Value = "2.73"
$c = $ws2.Cells.Item(3,7); $c.NumberFormat = "@"; $c.Value = "0.3969"
$ws2.Cells.Item(3,8).Value = 9
$ws2.Cells.Item(4,1).Value = 2
$c = $ws2.Cells.Item(4,2); $c.NumberFormat = "@"; $c.Value = "001239"
$ws2.Cells.Item(4,3).Value = "长盛国企改革主题灵活配置混合"
$c = $ws2.Cells.Item(4,4); $c.NumberFormat = "@"; $c.Value = "4.46"
$c = $ws2.Cells.Item(4,5); $c.NumberFormat = "@"; $c.Value = "90.97"
$c = $ws2.Cells.Item(4,6); $c.NumberFormat = "@"; $c.Value = "8.19"
$c = $ws2.Cells.Item(4,7); $c.NumberFormat = "@"; $c.Value = "0.3653"
$ws2.Cells.Item(4,8).Value = 3
$ws2.Cells.Item(5,1).Value = 3
$c = $ws2.Cells.Item(5,2); $c.NumberFormat = "@"; $c.Value = "010694"
$ws2.Cells.Item(5,3).Value = "万家内需增长一年持有期混合"
$c = $ws2.Cells.Item(5,4); $c.NumberFormat = "@"; $c.Value = "11.31"
$c = $ws2.Cells.Item(5,5); $c.NumberFormat = "@"; $c.Value = "95.03"
$c = $ws2.Cells.Item(5,6); $c.NumberFormat = "@"; $c.Value = "3.15"
$c = $ws2.Cells.Item(5,7); $c.NumberFormat = "@"; $c.Value = "0.3563"
$ws2.Cells.Item(5,8).Value = 10
$ws2.Cells.Item(6,1).Value = 4
$c = $ws2.Cells.Item(6,2); $c.NumberFormat = "@"; $c.Value = "005478"
$ws2.Cells.Item(6,3).Value = "长安鑫禧灵活配置混合C"
$c = $ws2.Cells.Item(6,4); $c.NumberFormat = "@"; $c.Value = "4.15"
$c = $ws2.Cells.Item(6,5); $c.NumberFormat = "@"; $c.Value = "91.79"
$c = $ws2.Cells.Item(6,6); $c.NumberFormat = "@"; $c.Value = "7.99"
$c = $ws2.Cells.Item(6,7); $c.NumberFormat = "@"; $c.Value = "0.3316"
$ws2.Cells.Item(6,8).Value = 3
$ws2.Cells.Item(7,1).Value = 5
$c = $ws2.Cells.Item(7,2); $c.NumberFormat = "@"; $c.Value = "501075"
$ws2.Cells.Item(7,3).Value = "万家科创主题灵活配置混合（LOF）A"
$c = $ws2.Cells.Item(7,4); $c.NumberFormat = "@"; $c.Value = "8.91"
$c = $ws2.Cells.Item(7,5); $c.NumberFormat = "@"; $c.Value = "94.55"
$c = $ws2.Cells.Item(7,6); $c.NumberFormat = "@"; $c.Value = "3.67"
$c = $ws2.Cells.Item(7,7); $c.NumberFormat = "@"; $c.Value = "0.3270"
$ws2.Cells.Item(7,8).Value = 7
$ws2.Cells.Item(8,1).Value = 6
$c = $ws2.Cells.Item(8,2); $c.NumberFormat = "@"; $c.Value = "005344"
$ws2.Cells.Item(8,3).Value = "长安裕盛灵活配置混合C"
$c = $ws2.Cells.Item(8,4); $c.NumberFormat = "@"; $c.Value = "3.87"
$c = $ws2.Cells.Item(8,5); $c.NumberFormat = "@"; $c.Value = "91.83"
$c = $ws2.Cells.Item(8,6); $c.NumberFormat = "@"; $c.Value = "7.79"
$c = $ws2.Cells.Item(8,7); $c.NumberFormat = "@"; $c.Value = "0.3015"
$ws2.Cells.Item(8,8).Value = 8
$ws2.Cells.Item(9,1).Value = 7
$c = $ws2.Cells.Item(9,2); $c.NumberFormat = "@"; $c.Value = "610004"
$ws2.Cells.Item(9,3).Value = "信澳中小盘混合"
$c = $ws2.Cells.Item(9,4); $c.NumberFormat = "@"; $c.Value = "5.35"
$c = $ws2.Cells.Item(9,5); $c.NumberFormat = "@"; $c.Value = "91.49"
$c = $ws2.Cells.Item(9,6); $c.NumberFormat = "@"; $c.Value = "5.32"
$c = $ws2.Cells.Item(9,7); $c.NumberFormat = "@"; $c.Value = "0.2846"
$ws2.Cells.Item(9,8).Value = 9
$ws2.Cells.Item(10,1).Value = 8
$c = $ws2.Cells.Item(10,2); $c.NumberFormat = "@"; $c.Value = "009859"
$ws2.Cells.Item(10,3).Value = "银华乐享混合A"
$c = $ws2.Cells.Item(10,4); $c.NumberFormat = "@"; $c.Value = "4.96"
$c = $ws2.Cells.Item(10,5); $c.NumberFormat = "@"; $c.Value = "94.27"
$c = $ws2.Cells.Item(10,6); $c.NumberFormat = "@"; $c.Value = "5.23"
$c = $ws2.Cells.Item(10,7); $c.NumberFormat = "@"; $c.Value = "0.2594"
$ws2.Cells.Item(10,8).Value = 6
$ws2.Cells.Item(11,1).Value = 9
$c = $ws2.Cells.Item(11,2); $c.NumberFormat = "@"; $c.Value = "005300"
$ws2.Cells.Item(11,3).Value = "万家成长优选灵活配置混合C"
$c = $ws2.Cells.Item(11,4); $c.NumberFormat = "@"; $c.Value = "8.08"
$c = $ws2.Cells.Item(11,5); $c.NumberFormat = "@"; $c.Value = "94.06"
$c = $ws2.Cells.Item(11,6); $c.NumberFormat = "@"; $c.Value = "2.94"
$c = $ws2.Cells.Item(11,7); $c.NumberFormat = "@"; $c.Value = "0.2376"
$ws2.Cells.Item(11,8).Value = 10
$ws2.Cells.Item(12,1).Value = 10
$c = $ws2.Cells.Item(12,2); $c.NumberFormat = "@"; $c.Value = "010611"
$ws2.Cells.Item(12,3).Value = "万家战略发展产业混合A"
$c = $ws2.Cells.Item(12,4); $c.NumberFormat = "@"; $c.Value = "6.69"
$c = $ws2.Cells.Item(12,5); $c.NumberFormat = "@"; $c.Value = "93.90"
$c = $ws2.Cells.Item(12,6); $c.NumberFormat = "@"; $c.Value = "3.39"
$c = $ws2.Cells.Item(12,7); $c.NumberFormat = "@"; $c.Value = "0.2268"
$ws2.Cells.Item(12,8).Value = 10
$ws2.Cells.Item(13,1).Value = 11
$c = $ws2.Cells.Item(13,2); $c.NumberFormat = "@"; $c.Value = "000800"
$ws2.Cells.Item(13,3).Value = "华商未来主题混合"
$c = $ws2.Cells.Item(13,4); $c.NumberFormat = "@"; $c.Value = "4.21"
$c = $ws2.Cells.Item(13,5); $c.NumberFormat = "@"; $c.Value = "72.31"
$c = $ws2.Cells.Item(13,6); $c.NumberFormat = "@"; $c.Value = "3.58"
$c = $ws2.Cells.Item(13,7); $c.NumberFormat = "@"; $c.Value = "0.1507"
$ws2.Cells.Item(13,8).Value = 8
$ws2.Cells.Item(14,1).Value = 12
$c = $ws2.Cells.Item(14,2); $c.NumberFormat = "@"; $c.Value = "010612"
$ws2.Cells.Item(14,3).Value = "万家战略发展产业混合C"
$c = $ws2.Cells.Item(14,4); $c.NumberFormat = "@"; $c.Value = "4.39"
$c = $ws2.Cells.Item(14,5); $c.NumberFormat = "@"; $c.Value = "93.90"
$c = $ws2.Cells.Item(14,6); $c.NumberFormat = "@"; $c.Value = "3.39"
$c = $ws2.Cells.Item(14,7); $c.NumberFormat = "@"; $c.Value = "0.1488"
$ws2.Cells.Item(14,8).Value = 10
$ws2.Cells.Item(15,1).Value = 13
$c = $ws2.Cells.Item(15,2); $c.NumberFormat = "@"; $c.Value = "161123"
$ws2.Cells.Item(15,3).Value = "易方达并购重组指数（LOF）"
$c = $ws2.Cells.Item(15,4); $c.NumberFormat = "@"; $c.Value = "4.34"
$c = $ws2.Cells.Item(15,5); $c.NumberFormat = "@"; $c.Value = "94.11"
$c = $ws2.Cells.Item(15,6); $c.NumberFormat = "@"; $c.Value = "2.62"
$c = $ws2.Cells.Item(15,7); $c.NumberFormat = "@"; $c.Value = "0.1137"
$ws2.Cells.Item(15,8).Value = 9
$ws2.Cells.Item(16,1).Value = 14
$c = $ws2.Cells.Item(16,2); $c.NumberFormat = "@"; $c.Value = "013495"
$ws2.Cells.Item(16,3).Value = "信澳产业优选一年持有混合A"
$c = $ws2.Cells.Item(16,4); $c.NumberFormat = "@"; $c.Value = "2.03"
$c = $ws2.Cells.Item(16,5); $c.NumberFormat = "@"; $c.Value = "79.40"
$c = $ws2.Cells.Item(16,6); $c.NumberFormat = "@"; $c.Value = "5.21"
$c = $ws2.Cells.Item(16,7); $c.NumberFormat = "@"; $c.Value = "0.1058"
$ws2.Cells.Item(16,8).Value = 5
$ws2.Cells.Item(17,1).Value = 15
$c = $ws2.Cells.Item(17,2); $c.NumberFormat = "@"; $c.Value = "001449"
$ws2.Cells.Item(17,3).Value = "华商双驱优选灵活配置混合"
$c = $ws2.Cells.Item(17,4); $c.NumberFormat = "@"; $c.Value = "2.41"
$c = $ws2.Cells.Item(17,5); $c.NumberFormat = "@"; $c.Value = "73.87"
$c = $ws2.Cells.Item(17,6); $c.NumberFormat = "@"; $c.Value = "4.30"
$c = $ws2.Cells.Item(17,7); $c.NumberFormat = "@"; $c.Value = "0.1036"
$ws2.Cells.Item(17,8).Value = 3
$ws2.Cells.Item(18,1).Value = 16
$c = $ws2.Cells.Item(18,2); $c.NumberFormat = "@"; $c.Value = "159625"
$ws2.Cells.Item(18,3).Value = "嘉实国证绿色电力ETF"
$c = $ws2.Cells.Item(18,4); $c.NumberFormat = "@"; $c.Value = "3.41"
$c = $ws2.Cells.Item(18,5); $c.NumberFormat = "@"; $c.Value = "98.77"
$c = $ws2.Cells.Item(18,6); $c.NumberFormat = "@"; $c.Value = "3.00"
$c = $ws2.Cells.Item(18,7); $c.NumberFormat = "@"; $c.Value = "0.1023"
$ws2.Cells.Item(18,8).Value = 9
$ws2.Cells.Item(19,1).Value = 17
$c = $ws2.Cells.Item(19,2); $c.NumberFormat = "@"; $c.Value = "005477"
$ws2.Cells.Item(19,3).Value = "长安鑫禧灵活配置混合A"
$c = $ws2.Cells.Item(19,4); $c.NumberFormat = "@"; $c.Value = "1.23"
$c = $ws2.Cells.Item(19,5); $c.NumberFormat = "@"; $c.Value = "91.79"
$c = $ws2.Cells.Item(19,6); $c.NumberFormat = "@"; $c.Value = "7.99"
$c = $ws2.Cells.Item(19,7); $c.NumberFormat = "@"; $c.Value = "0.0983"
$ws2.Cells.Item(19,8).Value = 3
$ws2.Cells.Item(20,1).Value = 18
$c = $ws2.Cells.Item(20,2); $c.NumberFormat = "@"; $c.Value = "006234"
$ws2.Cells.Item(20,3).Value = "万家汽车新趋势混合C"
$c = $ws2.Cells.Item(20,4); $c.NumberFormat = "@"; $c.Value = "2.23"
$c = $ws2.Cells.Item(20,5); $c.NumberFormat = "@"; $c.Value = "90.68"
$c = $ws2.Cells.Item(20,6); $c.NumberFormat = "@"; $c.Value = "3.50"
$c = $ws2.Cells.Item(20,7); $c.NumberFormat = "@"; $c.Value = "0.0780"
$ws2.Cells.Item(20,8).Value = 6
$ws2.Cells.Item(21,1).Value = 19
$c = $ws2.Cells.Item(21,2); $c.NumberFormat = "@"; $c.Value = "015687"
$ws2.Cells.Item(21,3).Value = "银华乐享混合C"
$c = $ws2.Cells.Item(21,4); $c.NumberFormat = "@"; $c.Value = "1.39"
$c = $ws2.Cells.Item(21,5); $c.NumberFormat = "@"; $c.Value = "94.27"
$c = $ws2.Cells.Item(21,6); $c.NumberFormat = "@"; $c.Value = "5.23"
$c = $ws2.Cells.Item(21,7); $c.NumberFormat = "@"; $c.Value = "0.0727"
$ws2.Cells.Item(21,8).Value = 6
$ws2.Cells.Item(22,1).Value = 20
$c = $ws2.Cells.Item(22,2); $c.NumberFormat = "@"; $c.Value = "006233"
$ws2.Cells.Item(22,3).Value = "万家汽车新趋势混合A"
$c = $ws2.Cells.Item(22,4); $c.NumberFormat = "@"; $c.Value = "1.93"
$c = $ws2.Cells.Item(22,5); $c.NumberFormat = "@"; $c.Value = "90.68"
$c = $ws2.Cells.Item(22,6); $c.NumberFormat = "@"; $c.Value = "3.50"
$c = $ws2.Cells.Item(22,7); $c.NumberFormat = "@"; $c.Value = "0.0676"
$ws2.Cells.Item(22,8).Value = 6
$ws2.Cells.Item(23,1).Value = 21
$c = $ws2.Cells.Item(23,2); $c.NumberFormat = "@"; $c.Value = "013326"
$ws2.Cells.Item(23,3).Value = "万家景气驱动混合A"
$c = $ws2.Cells.Item(23,4); $c.NumberFormat = "@"; $c.Value = "2.03"
$c = $ws2.Cells.Item(23,5); $c.NumberFormat = "@"; $c.Value = "92.93"
$c = $ws2.Cells.Item(23,6); $c.NumberFormat = "@"; $c.Value = "3.06"
$c = $ws2.Cells.Item(23,7); $c.NumberFormat = "@"; $c.Value = "0.0621"
$ws2.Cells.Item(23,8).Value = 9
$ws2.Cells.Item(24,1).Value = 22
$c = $ws2.Cells.Item(24,2); $c.NumberFormat = "@"; $c.Value = "005343"
$ws2.Cells.Item(24,3).Value = "长安裕盛灵活配置混合A"
$c = $ws2.Cells.Item(24,4); $c.NumberFormat = "@"; $c.Value = "0.62"
$c = $ws2.Cells.Item(24,5); $c.NumberFormat = "@"; $c.Value = "91.83"
$c = $ws2.Cells.Item(24,6); $c.NumberFormat = "@"; $c.Value = "7.79"
$c = $ws2.Cells.Item(24,7); $c.NumberFormat = "@"; $c.Value = "0.0483"
$ws2.Cells.Item(24,8).Value = 8
$ws2.Cells.Item(25,1).Value = 23
$c = $ws2.Cells.Item(25,2); $c.NumberFormat = "@"; $c.Value = "008602"
$ws2.Cells.Item(25,3).Value = "方正富邦新兴成长混合A"
$c = $ws2.Cells.Item(25,4); $c.NumberFormat = "@"; $c.Value = "1.23"
$c = $ws2.Cells.Item(25,5); $c.NumberFormat = "@"; $c.Value = "86.03"
$c = $ws2.Cells.Item(25,6); $c.NumberFormat = "@"; $c.Value = "3.67"
$c = $ws2.Cells.Item(25,7); $c.NumberFormat = "@"; $c.Value = "0.0451"
$ws2.Cells.Item(25,8).Value = 10
$ws2.Cells.Item(26,1).Value = 24
$c = $ws2.Cells.Item(26,2); $c.NumberFormat = "@"; $c.Value = "015384"
$ws2.Cells.Item(26,3).Value = "万家瑞隆混合C"
$c = $ws2.Cells.Item(26,4); $c.NumberFormat = "@"; $c.Value = "1.52"
$c = $ws2.Cells.Item(26,5); $c.NumberFormat = "@"; $c.Value = "93.96"
$c = $ws2.Cells.Item(26,6); $c.NumberFormat = "@"; $c.Value = "2.73"
$c = $ws2.Cells.Item(26,7); $c.NumberFormat = "@"; $c.Value = "0.0415"
$ws2.Cells.Item(26,8).Value = 9
$ws2.Cells.Item(27,1).Value = 25
$c = $ws2.Cells.Item(27,2); $c.NumberFormat = "@"; $c.Value = "002289"
$ws2.Cells.Item(27,3).Value = "华商改革创新股票A"
$c = $ws2.Cells.Item(27,4); $c.NumberFormat = "@"; $c.Value = "1.08"
$c = $ws2.Cells.Item(27,5); $c.NumberFormat = "@"; $c.Value = "79.89"
$c = $ws2.Cells.Item(27,6); $c.NumberFormat = "@"; $c.Value = "3.76"
$c = $ws2.Cells.Item(27,7); $c.NumberFormat = "@"; $c.Value = "0.0406"
$ws2.Cells.Item(27,8).Value = 9
$ws2.Cells.Item(28,1).Value = 26
$c = $ws2.Cells.Item(28,2); $c.NumberFormat = "@"; $c.Value = "001261"
$ws2.Cells.Item(28,3).Value = "中融新机遇灵活配置混合"
$c = $ws2.Cells.Item(28,4); $c.NumberFormat = "@"; $c.Value = "0.56"
$c = $ws2.Cells.Item(28,5); $c.NumberFormat = "@"; $c.Value = "93.44"
$c = $ws2.Cells.Item(28,6); $c.NumberFormat = "@"; $c.Value = "5.77"
$c = $ws2.Cells.Item(28,7); $c.NumberFormat = "@"; $c.Value = "0.0323"
$ws2.Cells.Item(28,8).Value = 8
$ws2.Cells.Item(29,1).Value = 27
$c = $ws2.Cells.Item(29,2); $c.NumberFormat = "@"; $c.Value = "000354"
$ws2.Cells.Item(29,3).Value = "长盛城镇化主题混合"
$c = $ws2.Cells.Item(29,4); $c.NumberFormat = "@"; $c.Value = "0.34"
$c = $ws2.Cells.Item(29,5); $c.NumberFormat = "@"; $c.Value = "92.91"
$c = $ws2.Cells.Item(29,6); $c.NumberFormat = "@"; $c.Value = "8.34"
$c = $ws2.Cells.Item(29,7); $c.NumberFormat = "@"; $c.Value = "0.0284"
$ws2.Cells.Item(29,8).Value = 1
$ws2.Cells.Item(30,1).Value = 28
$c = $ws2.Cells.Item(30,2); $c.NumberFormat = "@"; $c.Value = "010403"
$ws2.Cells.Item(30,3).Value = "华商景气优选混合"
$c = $ws2.Cells.Item(30,4); $c.NumberFormat = "@"; $c.Value = "0.61"
$c = $ws2.Cells.Item(30,5); $c.NumberFormat = "@"; $c.Value = "76.85"
$c = $ws2.Cells.Item(30,6); $c.NumberFormat = "@"; $c.Value = "4.35"
$c = $ws2.Cells.Item(30,7); $c.NumberFormat = "@"; $c.Value = "0.0265"
$ws2.Cells.Item(30,8).Value = 4
$ws2.Cells.Item(31,1).Value = 29
$c = $ws2.Cells.Item(31,2); $c.NumberFormat = "@"; $c.Value = "008491"
$ws2.Cells.Item(31,3).Value = "万家周期优势企业混合A"
$c = $ws2.Cells.Item(31,4); $c.NumberFormat = "@"; $c.Value = "0.61"
$c = $ws2.Cells.Item(31,5); $c.NumberFormat = "@"; $c.Value = "93.50"
$c = $ws2.Cells.Item(31,6); $c.NumberFormat = "@"; $c.Value = "3.10"
$c = $ws2.Cells.Item(31,7); $c.NumberFormat = "@"; $c.Value = "0.0189"
$ws2.Cells.Item(31,8).Value = 7
$ws2.Cells.Item(32,1).Value = 30
$c = $ws2.Cells.Item(32,2); $c.NumberFormat = "@"; $c.Value = "003704"
$ws2.Cells.Item(32,3).Value = "光大保德信事件驱动灵活配置混合"
$c = $ws2.Cells.Item(32,4); $c.NumberFormat = "@"; $c.Value = "1.69"
$c = $ws2.Cells.Item(32,5); $c.NumberFormat = "@"; $c.Value = "28.78"
$c = $ws2.Cells.Item(32,6); $c.NumberFormat = "@"; $c.Value = "1.11"
$c = $ws2.Cells.Item(32,7); $c.NumberFormat = "@"; $c.Value = "0.0188"
$ws2.Cells.Item(32,8).Value = 10
$ws2.Cells.Item(33,1).Value = 31
$c = $ws2.Cells.Item(33,2); $c.NumberFormat = "@"; $c.Value = "016052"
$ws2.Cells.Item(33,3).Value = "华商改革创新股票C"
$c = $ws2.Cells.Item(33,4); $c.NumberFormat = "@"; $c.Value = "0.32"
$c = $ws2.Cells.Item(33,5); $c.NumberFormat = "@"; $c.Value = "79.89"
$c = $ws2.Cells.Item(33,6); $c.NumberFormat = "@"; $c.Value = "3.76"
$c = $ws2.Cells.Item(33,7); $c.NumberFormat = "@"; $c.Value = "0.0120"
$ws2.Cells.Item(33,8).Value = 9
$ws2.Cells.Item(34,1).Value = 32
$c = $ws2.Cells.Item(34,2); $c.NumberFormat = "@"; $c.Value = "013327"
$ws2.Cells.Item(34,3).Value = "万家景气驱动混合C"
$c = $ws2.Cells.Item(34,4); $c.NumberFormat = "@"; $c.Value = "0.32"
$c = $ws2.Cells.Item(34,5); $c.NumberFormat = "@"; $c.Value = "92.93"
$c = $ws2.Cells.Item(34,6); $c.NumberFormat = "@"; $c.Value = "3.06"
$c = $ws2.Cells.Item(34,7); $c.NumberFormat = "@"; $c.Value = "0.0098"
$ws2.Cells.Item(34,8).Value = 9
$ws2.Cells.Item(35,1).Value = 33
$c = $ws2.Cells.Item(35,2); $c.NumberFormat = "@"; $c.Value = "013496"
$ws2.Cells.Item(35,3).Value = "信澳产业优选一年持有混合C"
$c = $ws2.Cells.Item(35,4); $c.NumberFormat = "@"; $c.Value = "0.18"
$c = $ws2.Cells.Item(35,5); $c.NumberFormat = "@"; $c.Value = "79.40"
$c = $ws2.Cells.Item(35,6); $c.NumberFormat = "@"; $c.Value = "5.21"
$c = $ws2.Cells.Item(35,7); $c.NumberFormat = "@"; $c.Value = "0.0094"
$ws2.Cells.Item(35,8).Value = 5
$ws2.Cells.Item(36,1).Value = 34
$c = $ws2.Cells.Item(36,2); $c.NumberFormat = "@"; $c.Value = "008492"
$ws2.Cells.Item(36,3).Value = "万家周期优势企业混合C"
$c = $ws2.Cells.Item(36,4); $c.NumberFormat = "@"; $c.Value = "0.14"
$c = $ws2.Cells.Item(36,5); $c.NumberFormat = "@"; $c.Value = "93.50"
$c = $ws2.Cells.Item(36,6); $c.NumberFormat = "@"; $c.Value = "3.10"
$c = $ws2.Cells.Item(36,7); $c.NumberFormat = "@"; $c.Value = "0.0043"
$ws2.Cells.Item(36,8).Value = 7
$ws2.Cells.Item(37,1).Value = 35
$c = $ws2.Cells.Item(37,2); $c.NumberFormat = "@"; $c.Value = "007501"
$ws2.Cells.Item(37,3).Value = "万家科创主题灵活配置混合（LOF）C"
$c = $ws2.Cells.Item(37,4); $c.NumberFormat = "@"; $c.Value = "0.09"
$c = $ws2.Cells.Item(37,5); $c.NumberFormat = "@"; $c.Value = "94.55"
$c = $ws2.Cells.Item(37,6); $c.NumberFormat = "@"; $c.Value = "3.67"
$c = $ws2.Cells.Item(37,7); $c.NumberFormat = "@"; $c.Value = "0.0033"
$ws2.Cells.Item(37,8).Value = 7
$ws2.Cells.Item(38,1).Value = 36
$c = $ws2.Cells.Item(38,2); $c.NumberFormat = "@"; $c.Value = "001899"
$ws2.Cells.Item(38,3).Value = "东海中证社会发展安全产业主题指数"
$c = $ws2.Cells.Item(38,4); $c.NumberFormat = "@"; $c.Value = "0.19"
$c = $ws2.Cells.Item(38,5); $c.NumberFormat = "@"; $c.Value = "90.94"
$c = $ws2.Cells.Item(38,6); $c.NumberFormat = "@"; $c.Value = "1.53"
$c = $ws2.Cells.Item(38,7); $c.NumberFormat = "@"; $c.Value = "0.0029"
$ws2.Cells.Item(38,8).Value = 8
$ws2.Cells.Item(39,1).Value = 37
$c = $ws2.Cells.Item(39,2); $c.NumberFormat = "@"; $c.Value = "011987"
$ws2.Cells.Item(39,3).Value = "财通资管智选核心回报6个月持有期混合A"
$c = $ws2.Cells.Item(39,4); $c.NumberFormat = "@"; $c.Value = "0.12"
$c = $ws2.Cells.Item(39,5); $c.NumberFormat = "@"; $c.Value = "39.46"
$c = $ws2.Cells.Item(39,6); $c.NumberFormat = "@"; $c.Value = "1.17"
$c = $ws2.Cells.Item(39,7); $c.NumberFormat = "@"; $c.Value = "0.0014"
$ws2.Cells.Item(39,8).Value = 5
$ws2.Cells.Item(40,1).Value = 38
$c = $ws2.Cells.Item(40,2); $c.NumberFormat = "@"; $c.Value = "008603"
$ws2.Cells.Item(40,3).Value = "方正富邦新兴成长混合C"
$c = $ws2.Cells.Item(40,4); $c.NumberFormat = "@"; $c.Value = "0.03"
$c = $ws2.Cells.Item(40,5); $c.NumberFormat = "@"; $c.Value = "86.03"
$c = $ws2.Cells.Item(40,6); $c.NumberFormat = "@"; $c.Value = "3.67"
$c = $ws2.Cells.Item(40,7); $c.NumberFormat = "@"; $c.Value = "0.0011"
$ws2.Cells.Item(40,8).Value = 10
$ws2.Cells.Item(41,1).Value = 39
$c = $ws2.Cells.Item(41,2); $c.NumberFormat = "@"; $c.Value = "011988"
$ws2.Cells.Item(41,3).Value = "财通资管智选核心回报6个月持有期混合C"
$c = $ws2.Cells.Item(41,4); $c.NumberFormat = "@"; $c.Value = "0.01"
$c = $ws2.Cells.Item(41,5); $c.NumberFormat = "@"; $c.Value = "39.46"
$c = $ws2.Cells.Item(41,6); $c.NumberFormat = "@"; $c.Value = "1.17"
$c = $ws2.Cells.Item(41,7); $c.NumberFormat = "@"; $c.Value = "0.0001"
$ws2.Cells.Item(41,8).Value = 5

# =========================================================
# Step 2: update the "总计" (summary) sheet -- insert a new
# row for 2022-Q3 at the top of the data and push every other
# quarter row down by one (re-numbering column A as we go).
# =========================================================
$ws1 = $wb.Worksheets.Item("总计")

# clone the numbering-column style down onto the new row 6
$ws1.Cells.Item(5,1).Copy()
$ws1.Cells.Item(6,1).PasteSpecial(-4122)

# fill bottom-up so each row's source values are read before they
# get overwritten by the row above
$ws1.Cells.Item(6,1).Value = 4
$ws1.Cells.Item(6,2).Value = "2021-Q3"
$ws1.Cells.Item(6,3).Value = 13
$ws1.Cells.Item(6,4).Value = 2.94

$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = "2021-Q4"
$ws1.Cells.Item(5,3).Value = 8
$ws1.Cells.Item(5,4).Value = 1.11

$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "2022-Q1"
$ws1.Cells.Item(4,3).Value = 5
$ws1.Cells.Item(4,4).Value = 0.96

$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q2"
$ws1.Cells.Item(3,3).Value = 8
$ws1.Cells.Item(3,4).Value = 0.92

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q3"
$ws1.Cells.Item(2,3).Value = 40
$ws1.Cells.Item(2,4).Value = 4.96

Write-Host "2022-Q3 sheet inserted and 总计 summary updated"
